$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23515522480011
$ws.Range("B1").Value = 2.548962593078613
$ws.Range("C1").Value = 9.146985054016113
$ws.Range("D1").Value = 2.035535097122192
$ws.Range("E1").Value = 1.171088218688965
